# Generate Report for Handoff
#
# Inserts a new "handoff" row (for file 0d1ca4e8-563f-4906-8bcc-a3977a07398f)
# above the existing a0dfcb73-3710-42e5-b5b0-373815c853ed row on every
# worksheet (Overview, zh-cn, de-de) of the localization-status report.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Remove existing hyperlinks before shifting rows around - the underlying
# hyperlink ranges do not follow inserted/moved cells automatically.
$ws.Hyperlinks.Delete()

# Push the existing data row down and make room for the new row.
$ws.Rows(2).Insert()

$ws.Cells.Item(2,1).Value = "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md"
$ws.Cells.Item(2,2).Value = "Ready for handoff"
$ws.Cells.Item(2,3).Value = "Ready for handoff"
$ws.Cells.Item(2,4).Value = "2016-29-19 06:29:03"

# Row 3 already holds the original a0dfcb73... values (shifted down by the
# insert above); nothing else needs to change there.

$ws.Hyperlinks.Add($ws.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/0d1ca4e8-563f-4906-8bcc-a3977a07398f.md", "", "", "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/a0dfcb73-3710-42e5-b5b0-373815c853ed.md", "", "", "a0dfcb73-3710-42e5-b5b0-373815c853ed.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()
$ws.Rows(2).Insert()

$ws.Cells.Item(2,1).Value = "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md"
$ws.Cells.Item(2,2).Value = ".md"
$ws.Cells.Item(2,3).Value = "Ready for handoff"
$ws.Cells.Item(2,4).Value = "0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.zh-cn.xlf"
$ws.Cells.Item(2,5).Value = "2016-03-19 06:28:59"
$ws.Cells.Item(2,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Cells.Item(2,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(2,9).Value = "Include"

# Row 3 already holds the original a0dfcb73... values (shifted down).

$ws.Hyperlinks.Add($ws.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/0d1ca4e8-563f-4906-8bcc-a3977a07398f.md", "", "", "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(2,2), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/0d1ca4e8-563f-4906-8bcc-a3977a07398f.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(2,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b40431e5f8e51d6cdae64b193740bad9d014da95/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.zh-cn.xlf", "", "", "0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/a0dfcb73-3710-42e5-b5b0-373815c853ed.md", "", "", "a0dfcb73-3710-42e5-b5b0-373815c853ed.md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3,2), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/a0dfcb73-3710-42e5-b5b0-373815c853ed.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b40431e5f8e51d6cdae64b193740bad9d014da95/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a0dfcb73-3710-42e5-b5b0-373815c853ed.7c23583881e90434debdf5bd12e534d97478fab2.zh-cn.xlf", "", "", "a0dfcb73-3710-42e5-b5b0-373815c853ed.7c23583881e90434debdf5bd12e534d97478fab2.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()
$ws.Rows(2).Insert()

$ws.Cells.Item(2,1).Value = "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md"
$ws.Cells.Item(2,2).Value = ".md"
$ws.Cells.Item(2,3).Value = "Ready for handoff"
$ws.Cells.Item(2,4).Value = "0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.de-de.xlf"
$ws.Cells.Item(2,5).Value = "2016-03-19 06:29:03"
$ws.Cells.Item(2,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Cells.Item(2,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(2,9).Value = "Include"

# Row 3 already holds the original a0dfcb73... values (shifted down).

$ws.Hyperlinks.Add($ws.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/0d1ca4e8-563f-4906-8bcc-a3977a07398f.md", "", "", "0d1ca4e8-563f-4906-8bcc-a3977a07398f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(2,2), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/0d1ca4e8-563f-4906-8bcc-a3977a07398f.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(2,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef0aef0432d93019bea41c1cc46a73929fdaa4fc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.de-de.xlf", "", "", "0d1ca4e8-563f-4906-8bcc-a3977a07398f.ccb5949cadfba9ae28124f850d36e8217cf49b07.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/a0dfcb73-3710-42e5-b5b0-373815c853ed.md", "", "", "a0dfcb73-3710-42e5-b5b0-373815c853ed.md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3,2), "https://github.com/OpenLocalizationTest/oltest/blob/7d71ef2bcd378c0e89219ad9d8df1faac5e6ec0a/e2e/a0dfcb73-3710-42e5-b5b0-373815c853ed.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef0aef0432d93019bea41c1cc46a73929fdaa4fc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a0dfcb73-3710-42e5-b5b0-373815c853ed.7c23583881e90434debdf5bd12e534d97478fab2.de-de.xlf", "", "", "a0dfcb73-3710-42e5-b5b0-373815c853ed.7c23583881e90434debdf5bd12e534d97478fab2.de-de.xlf") | Out-Null

$wb.Save()
